# Add the "attendance" sheet (Punch In/Out time tracking) after "timeSheet".
# Mirrors the author's commit "added code for attendance".

$wb = $excel.ActiveWorkbook
$timeSheet = $wb.Worksheets.Item(1)

# --- Shrink the width of timeSheet's "execution status" column slightly ---
$timeSheet.Columns.Item(6).ColumnWidth = 12.6

# --- Create the new worksheet right after timeSheet ---
$attendance = $wb.Worksheets.Add($null, $timeSheet)
$attendance.Name = "attendance"

# Header row
$attendance.Range("A1").Value = "Testcase"
$attendance.Range("B1").Value = "Run"
$attendance.Range("C1").Value = "hours"
$attendance.Range("D1").Value = "'Punch In Time"
$attendance.Range("E1").Value = "'Punch Out Time"
$attendance.Range("F1").Value = "execution status"

# Data row 2 (shared-string order: Punch In/Out headers, then the two punch
# times, then the elapsed-hours text, matching the author's entry order).
# Value is written before NumberFormat so the text-with-apostrophe (quote
# prefix) style is the only one ever created for each cell.
$attendance.Range("A2").Value = 1
$attendance.Range("B2").Value = "yes"
$attendance.Range("D2").Value = "'09:00 AM"
$attendance.Range("D2").NumberFormat = "h:mm AM/PM"
$attendance.Range("E2").Value = "'06:00 PM"
$attendance.Range("E2").NumberFormat = "h:mm AM/PM"
$attendance.Range("C2").Value = "'9:00"
$attendance.Range("C2").NumberFormat = "h:mm"
$attendance.Range("F2").Value = "done"

# Remaining rows
$attendance.Range("A3").Value = 2
$attendance.Range("B3").Value = "no"
$attendance.Range("A4").Value = 3
$attendance.Range("B4").Value = "no"
$attendance.Range("A5").Value = 4
$attendance.Range("B5").Value = "no"

# Column widths / page setup for the new sheet
$attendance.Columns.Item(4).ColumnWidth = 11.85
$attendance.PageSetup.PaperSize = 9
$attendance.PageSetup.Orientation = 1

# Selection: active cell E2 on the attendance tab (matches author's last click)
$attendance.Range("E2").Select()

# timeSheet keeps a full-range selection once it is no longer the active tab
$timeSheet.Range("A1:F5").Select()
$attendance.Activate()
